$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 1209.9683
$ws.Range("I15").Value = 1209.9683
$ws.Range("K15").Value = 3629.9049
$ws.Range("M15").Value = -3460.9049
# Row 21
$ws.Range("H21").Value = 34999.25
$ws.Range("I21").Value = 29999
$ws.Range("J21").Value = 50000
$ws.Range("K21").Value = 29999
$ws.Range("L21").Value = 50000
$ws.Range("M21").Value = -29531
$ws.Range("N21").Value = -50936
# Row 23
$ws.Range("H23").Value = 34999.25
$ws.Range("I23").Value = 29999
$ws.Range("J23").Value = 50000
$ws.Range("K23").Value = 29999
$ws.Range("L23").Value = 50000
$ws.Range("M23").Value = -29765
$ws.Range("N23").Value = -50468
# Row 29
$ws.Range("H29").Value = 460.6
$ws.Range("J29").Value = 2004
$ws.Range("L29").Value = 6012
$ws.Range("N29").Value = -6574
# Row 38
$ws.Range("H38").Value = 1348.8636
$ws.Range("I38").Value = 57.5
$ws.Range("J38").Value = 2425
$ws.Range("K38").Value = 172.5
$ws.Range("L38").Value = 7275
$ws.Range("M38").Value = 199.5
$ws.Range("N38").Value = -8019
# Row 43
$ws.Range("H43").Value = 1210
$ws.Range("I43").Value = 600
$ws.Range("J43").Value = 1820
$ws.Range("K43").Value = 600
$ws.Range("L43").Value = 1820
$ws.Range("M43").Value = -531
$ws.Range("N43").Value = -1958
# Row 51
$ws.Range("H51").Value = 2135.5454
$ws.Range("I51").Value = 1527.2858
$ws.Range("J51").Value = 3200
$ws.Range("K51").Value = 1527.2858
$ws.Range("L51").Value = 3200
$ws.Range("M51").Value = -1043.2858
$ws.Range("N51").Value = -4168
# Row 64
$ws.Range("H64").Value = 3287.3
$ws.Range("I64").Value = 3059.9167
$ws.Range("J64").Value = 3438.889
$ws.Range("K64").Value = 3059.9167
$ws.Range("L64").Value = 3438.889
$ws.Range("M64").Value = -2811.9167
$ws.Range("N64").Value = -3934.889
# Row 67
$ws.Range("H67").Value = 3287.3
$ws.Range("I67").Value = 3059.9167
$ws.Range("J67").Value = 3438.889
$ws.Range("K67").Value = 3059.9167
$ws.Range("L67").Value = 3438.889
$ws.Range("M67").Value = -2201.9167
$ws.Range("N67").Value = -5154.889
# Row 94
$ws.Range("H94").Value = 6201.1113
$ws.Range("I94").Value = 6201.1113
$ws.Range("K94").Value = 6201.1113
$ws.Range("M94").Value = -5750.1113
# Row 100
$ws.Range("H100").Value = 3483.75
$ws.Range("I100").Value = 3000
$ws.Range("J100").Value = 3580.5
$ws.Range("K100").Value = 3000
$ws.Range("L100").Value = 3580.5
$ws.Range("M100").Value = -2459
$ws.Range("N100").Value = -4662.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 21
$ws.Range("H21").Value = 23110
$ws.Range("J21").Value = 18000
$ws.Range("L21").Value = 18000
$ws.Range("N21").Value = -18748
# Row 32
$ws.Range("H32").Value = 13010.8125
$ws.Range("I32").Value = 13202.532
$ws.Range("J32").Value = 4000
$ws.Range("K32").Value = 13202.532
$ws.Range("L32").Value = 4000
$ws.Range("M32").Value = -12915.532
$ws.Range("N32").Value = -4574
# Row 60
$ws.Range("H60").Value = 28025.5
$ws.Range("I60").Value = 28025.5
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 28025.5
$ws.Range("L60").Value = 0
$ws.Range("M60").Value = -27292.5
$ws.Range("N60").ClearContents()
# Row 70
$ws.Range("H70").Value = 50000
$ws.Range("J70").Value = 50000
$ws.Range("L70").Value = 50000
$ws.Range("N70").Value = -50540
# Row 73
$ws.Range("H73").Value = 50000
$ws.Range("J73").Value = 50000
$ws.Range("L73").Value = 50000
$ws.Range("N73").Value = -51872
# Row 74
$ws.Range("H74").Value = 1360.0667
$ws.Range("I74").Value = 1060.6
$ws.Range("J74").Value = 1959
$ws.Range("K74").Value = 1060.6
$ws.Range("L74").Value = 1959
$ws.Range("M74").Value = -186.5999999999999
$ws.Range("N74").Value = -3707
# Row 77
$ws.Range("H77").Value = 1360.0667
$ws.Range("I77").Value = 1060.6
$ws.Range("J77").Value = 1959
$ws.Range("K77").Value = 5303
$ws.Range("L77").Value = 9795
$ws.Range("M77").Value = -935
$ws.Range("N77").Value = -18531
# Row 140
$ws.Range("H140").Value = 52416.09
$ws.Range("J140").Value = 52416.09
$ws.Range("L140").Value = 52416.09
$ws.Range("N140").Value = -62776.09

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 59
$ws.Range("H59").Value = 62826.668
$ws.Range("J59").Value = 62826.668
$ws.Range("L59").Value = 62826.668
$ws.Range("N59").Value = -64520.668
# Row 95
$ws.Range("H95").Value = 22455
$ws.Range("J95").Value = 22455
$ws.Range("L95").Value = 22455
$ws.Range("N95").Value = -27947

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 37
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()
# Row 56
$ws.Range("H56").Value = 23000
$ws.Range("I56").Value = 29500
$ws.Range("J56").Value = 10000
$ws.Range("K56").Value = 29500
$ws.Range("L56").Value = 10000
$ws.Range("M56").Value = -28655
$ws.Range("N56").Value = -11690
# Row 99
$ws.Range("H99").Value = 2867.5715
$ws.Range("I99").Value = 2918.1875
$ws.Range("J99").Value = 2705.6
$ws.Range("K99").Value = 2918.1875
$ws.Range("L99").Value = 2705.6
$ws.Range("M99").Value = -1420.1875
$ws.Range("N99").Value = -5701.6
# Row 126
$ws.Range("H126").Value = 2867.5715
$ws.Range("I126").Value = 2918.1875
$ws.Range("J126").Value = 2705.6
$ws.Range("K126").Value = 8754.5625
$ws.Range("L126").Value = 8116.799999999999
$ws.Range("M126").Value = -6284.5625
$ws.Range("N126").Value = -13056.8

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 113
$ws.Range("H113").Value = 417315.6
$ws.Range("I113").Value = 833982.06
$ws.Range("J113").Value = 649.0833
$ws.Range("K113").Value = 2501946.18
$ws.Range("L113").Value = 1947.2499
$ws.Range("M113").Value = -2499776.18
$ws.Range("N113").Value = -6287.2499
# Row 120
$ws.Range("H120").Value = 8850.533
$ws.Range("I120").Value = 5939.778
$ws.Range("J120").Value = 13216.667
$ws.Range("K120").Value = 17819.334
$ws.Range("L120").Value = 39650.001
$ws.Range("M120").Value = -12981.334
$ws.Range("N120").Value = -49326.001
# Row 131
$ws.Range("H131").Value = 14085908
$ws.Range("I131").Value = 370
$ws.Range("J131").Value = 14707329
$ws.Range("K131").Value = 1110
$ws.Range("L131").Value = 44121987
$ws.Range("M131").Value = 3930
$ws.Range("N131").Value = -44132067

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 1496.6666
$ws.Range("I16").Value = 1496.6666
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1496.6666
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -1326.6666
$ws.Range("N16").ClearContents()
# Row 68
$ws.Range("H68").Value = 2048.5
$ws.Range("J68").Value = 2324.3635
$ws.Range("L68").Value = 2324.3635
$ws.Range("N68").Value = -3822.3635
# Row 71
$ws.Range("H71").Value = 2048.5
$ws.Range("J71").Value = 2324.3635
$ws.Range("L71").Value = 11621.8175
$ws.Range("N71").Value = -19109.8175
# Row 93
$ws.Range("H93").Value = 1222.2222
# Row 94
$ws.Range("H94").Value = 35000
$ws.Range("J94").Value = 35000
$ws.Range("L94").Value = 35000
$ws.Range("N94").Value = -36352
# Row 132
$ws.Range("H132").Value = 4260.9463
$ws.Range("I132").Value = 3744.9429
$ws.Range("J132").Value = 5120.952
$ws.Range("K132").Value = 11234.8287
$ws.Range("L132").Value = 15362.856
$ws.Range("M132").Value = -8704.8287
$ws.Range("N132").Value = -20422.856

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 1485.7142
$ws.Range("J96").Value = 1600
$ws.Range("L96").Value = 1600
$ws.Range("N96").Value = -4346
# Row 105
$ws.Range("H105").Value = 28333.334
$ws.Range("J105").Value = 28333.334
$ws.Range("L105").Value = 28333.334
$ws.Range("N105").Value = -35321.334
# Row 123
$ws.Range("H123").Value = 26668
$ws.Range("J123").Value = 26668
$ws.Range("L123").Value = 26668
$ws.Range("N123").Value = -36468
